$wb = $excel.ActiveWorkbook

# Map of worksheet name -> new "gewijzigd" (modified) timestamp serial value in column B11
$timestamps = @{
    "IK106" = 45096.4670691518
    "Q100"  = 45096.46565826389
    "Q200"  = 45096.46573024306
    "Q300"  = 45096.46580151621
    "Q400"  = 45096.46589184028
    "Q500"  = 45096.46598178241
    "Q600"  = 45096.46604519676
    "P100"  = 45096.4661228125
    "P200"  = 45096.46621456018
    "P300"  = 45096.46626951389
    "P400"  = 45096.46633971065
    "P500"  = 45096.46642256944
    "P600"  = 45096.46650778935
    "IK91"  = 45096.46657356482
    "IK92"  = 45096.46661386574
    "IK93"  = 45096.46666118056
    "IK94"  = 45096.46671474537
    "IK95"  = 45096.46676322917
    "IK96"  = 45096.46680692129
    "IK101" = 45096.46684737269
    "IK102" = 45096.46689071759
    "IK103" = 45096.46693547453
    "IK104" = 45096.46698012731
    "IK105" = 45096.46702700231
}

foreach ($name in $timestamps.Keys) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("B11").Value = $timestamps[$name]
}

# Sheet P100 got additional refitted model coefficients
$wsP100 = $wb.Worksheets.Item("P100")
$wsP100.Range("B2").Value = -0.008696875167344879
$wsP100.Range("B4").Value = -0.0000000001
$wsP100.Range("B6").Value = 4.29653460926825
$wsP100.Range("B7").Value = 170.3685582634114
$wsP100.Range("B10").Value = 0.6199068241255263
